$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 58328.055
$ws.Cells.Item(86, 9).Value = 202581.2
$ws.Cells.Item(86, 10).Value = 2846.077
$ws.Cells.Item(86, 11).Value = 202581.2
$ws.Cells.Item(86, 12).Value = 2846.077
$ws.Cells.Item(86, 13).Value = -201458.2
$ws.Cells.Item(86, 14).Value = -5092.077

$ws.Cells.Item(88, 8).Value = 1248.4445
$ws.Cells.Item(88, 9).Value = 845.3333
$ws.Cells.Item(88, 10).Value = 1450
$ws.Cells.Item(88, 11).Value = 845.3333
$ws.Cells.Item(88, 12).Value = 1450
$ws.Cells.Item(88, 13).Value = -439.3333
$ws.Cells.Item(88, 14).Value = -2262

$ws.Cells.Item(89, 8).Value = 58328.055
$ws.Cells.Item(89, 9).Value = 202581.2
$ws.Cells.Item(89, 10).Value = 2846.077
$ws.Cells.Item(89, 11).Value = 1012906
$ws.Cells.Item(89, 12).Value = 14230.385
$ws.Cells.Item(89, 13).Value = -1007290
$ws.Cells.Item(89, 14).Value = -25462.385

$ws.Cells.Item(91, 8).Value = 1248.4445
$ws.Cells.Item(91, 9).Value = 845.3333
$ws.Cells.Item(91, 10).Value = 1450
$ws.Cells.Item(91, 11).Value = 845.3333
$ws.Cells.Item(91, 12).Value = 1450
$ws.Cells.Item(91, 13).Value = 558.6667
$ws.Cells.Item(91, 14).Value = -4258

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6249.7646
$ws.Cells.Item(132, 9).Value = 2135.1
$ws.Cells.Item(132, 10).Value = 7964.2085
$ws.Cells.Item(132, 11).Value = 6405.299999999999
$ws.Cells.Item(132, 12).Value = 23892.6255
$ws.Cells.Item(132, 13).Value = -3875.299999999999
$ws.Cells.Item(132, 14).Value = -28952.6255

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 55780
$ws.Cells.Item(52, 10).Value = 55780
$ws.Cells.Item(52, 12).Value = 55780
$ws.Cells.Item(52, 14).Value = -56306

$ws.Cells.Item(86, 8).Value = 2266.4
$ws.Cells.Item(86, 9).Value = 2299.7036
$ws.Cells.Item(86, 10).Value = 1966.6666
$ws.Cells.Item(86, 11).Value = 2299.7036
$ws.Cells.Item(86, 12).Value = 1966.6666
$ws.Cells.Item(86, 13).Value = -1176.7036
$ws.Cells.Item(86, 14).Value = -4212.6666

$ws.Cells.Item(89, 8).Value = 2266.4
$ws.Cells.Item(89, 9).Value = 2299.7036
$ws.Cells.Item(89, 10).Value = 1966.6666
$ws.Cells.Item(89, 11).Value = 11498.518
$ws.Cells.Item(89, 12).Value = 9833.333000000001
$ws.Cells.Item(89, 13).Value = -5882.518
$ws.Cells.Item(89, 14).Value = -21065.333

$ws.Cells.Item(121, 8).Value = 55780
$ws.Cells.Item(121, 10).Value = 55780
$ws.Cells.Item(121, 12).Value = 55780
$ws.Cells.Item(121, 14).Value = -59274

$ws.Cells.Item(134, 8).Value = 4205.8823
$ws.Cells.Item(134, 9).Value = 4030.303
$ws.Cells.Item(134, 10).Value = 10000
$ws.Cells.Item(134, 11).Value = 12090.909
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 13).Value = -9555.909
$ws.Cells.Item(134, 14).Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 51402.2
$ws.Cells.Item(10, 9).Value = 499.5
$ws.Cells.Item(10, 10).Value = 64127.875
$ws.Cells.Item(10, 11).Value = 499.5
$ws.Cells.Item(10, 12).Value = 64127.875
$ws.Cells.Item(10, 13).Value = -360.5
$ws.Cells.Item(10, 14).Value = -64405.875

$ws.Cells.Item(31, 8).Value = 1760.1666
$ws.Cells.Item(31, 9).Value = 1181.8383
$ws.Cells.Item(31, 10).Value = 5692.8
$ws.Cells.Item(31, 11).Value = 1181.8383
$ws.Cells.Item(31, 12).Value = 5692.8
$ws.Cells.Item(31, 13).Value = -886.8382999999999
$ws.Cells.Item(31, 14).Value = -6282.8

$ws.Cells.Item(34, 8).Value = 1760.1666
$ws.Cells.Item(34, 9).Value = 1181.8383
$ws.Cells.Item(34, 10).Value = 5692.8
$ws.Cells.Item(34, 11).Value = 1181.8383
$ws.Cells.Item(34, 12).Value = 5692.8
$ws.Cells.Item(34, 13).Value = -979.8382999999999
$ws.Cells.Item(34, 14).Value = -6096.8

$ws.Cells.Item(58, 8).Value = 4136648.5
$ws.Cells.Item(58, 9).Value = 11367110
$ws.Cells.Item(58, 10).Value = 4956.357
$ws.Cells.Item(58, 11).Value = 11367110
$ws.Cells.Item(58, 12).Value = 4956.357
$ws.Cells.Item(58, 13).Value = -11366907
$ws.Cells.Item(58, 14).Value = -5362.357

$ws.Cells.Item(132, 8).Value = 2672.913
$ws.Cells.Item(132, 9).Value = 2273.9
$ws.Cells.Item(132, 10).Value = 5333
$ws.Cells.Item(132, 11).Value = 6821.700000000001
$ws.Cells.Item(132, 12).Value = 15999
$ws.Cells.Item(132, 13).Value = -4291.700000000001
$ws.Cells.Item(132, 14).Value = -21059

$ws.Cells.Item(134, 8).Value = 4381.6216
$ws.Cells.Item(134, 9).Value = 3294.3572
$ws.Cells.Item(134, 10).Value = 5043.4346
$ws.Cells.Item(134, 11).Value = 9883.071599999999
$ws.Cells.Item(134, 12).Value = 15130.3038
$ws.Cells.Item(134, 13).Value = -7348.071599999999
$ws.Cells.Item(134, 14).Value = -20200.3038

$ws.Cells.Item(136, 8).Value = 4136648.5
$ws.Cells.Item(136, 9).Value = 11367110
$ws.Cells.Item(136, 10).Value = 4956.357
$ws.Cells.Item(136, 11).Value = 34101330
$ws.Cells.Item(136, 12).Value = 14869.071
$ws.Cells.Item(136, 13).Value = -34098780
$ws.Cells.Item(136, 14).Value = -19969.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 27607.111
$ws.Cells.Item(131, 9).Value = 1897.5
$ws.Cells.Item(131, 10).Value = 37495.42
$ws.Cells.Item(131, 11).Value = 5692.5
$ws.Cells.Item(131, 12).Value = 112486.26
$ws.Cells.Item(131, 13).Value = -652.5
$ws.Cells.Item(131, 14).Value = -122566.26

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 12750
$ws.Cells.Item(31, 9).Value = 500
$ws.Cells.Item(31, 10).Value = 25000
$ws.Cells.Item(31, 11).Value = 500
$ws.Cells.Item(31, 12).Value = 25000
$ws.Cells.Item(31, 13).Value = -208
$ws.Cells.Item(31, 14).Value = -25584

$ws.Cells.Item(37, 8).Value = 12750
$ws.Cells.Item(37, 9).Value = 500
$ws.Cells.Item(37, 10).Value = 25000
$ws.Cells.Item(37, 11).Value = 500
$ws.Cells.Item(37, 12).Value = 25000
$ws.Cells.Item(37, 13).Value = -223
$ws.Cells.Item(37, 14).Value = -25554

$ws.Cells.Item(102, 8).Value = 6098.75
$ws.Cells.Item(102, 9).Value = 6200
$ws.Cells.Item(102, 10).Value = 5930
$ws.Cells.Item(102, 11).Value = 6200
$ws.Cells.Item(102, 12).Value = 5930
$ws.Cells.Item(102, 13).Value = -4578
$ws.Cells.Item(102, 14).Value = -9174

$ws.Cells.Item(123, 8).Value = 27771.8
$ws.Cells.Item(123, 10).Value = 27771.8
$ws.Cells.Item(123, 12).Value = 27771.8
$ws.Cells.Item(123, 14).Value = -32671.8

$ws.Cells.Item(126, 8).Value = 3666.3333
$ws.Cells.Item(126, 9).Value = 2000
$ws.Cells.Item(126, 10).Value = 4777.222
$ws.Cells.Item(126, 11).Value = 6000
$ws.Cells.Item(126, 12).Value = 14331.666
$ws.Cells.Item(126, 13).Value = -3530
$ws.Cells.Item(126, 14).Value = -19271.666

$ws.Cells.Item(132, 8).Value = 5472.364
$ws.Cells.Item(132, 9).Value = 5885.857
$ws.Cells.Item(132, 10).Value = 4748.75
$ws.Cells.Item(132, 11).Value = 17657.571
$ws.Cells.Item(132, 12).Value = 14246.25
$ws.Cells.Item(132, 13).Value = -15127.571
$ws.Cells.Item(132, 14).Value = -19306.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6908.3335
$ws.Cells.Item(7, 9).Value = 4900
$ws.Cells.Item(7, 10).Value = 8342.857
$ws.Cells.Item(7, 11).Value = 4900
$ws.Cells.Item(7, 12).Value = 8342.857
$ws.Cells.Item(7, 13).Value = -4788
$ws.Cells.Item(7, 14).Value = -8566.857

$ws.Cells.Item(40, 8).Value = 4233.952
$ws.Cells.Item(40, 9).Value = 4059.2942
$ws.Cells.Item(40, 11).Value = 4059.2942
$ws.Cells.Item(40, 13).Value = -3923.2942

$ws.Cells.Item(68, 8).Value = 1620
$ws.Cells.Item(68, 9).Value = 1550
$ws.Cells.Item(68, 11).Value = 1550
$ws.Cells.Item(68, 13).Value = -801

$ws.Cells.Item(71, 8).Value = 1620
$ws.Cells.Item(71, 9).Value = 1550
$ws.Cells.Item(71, 11).Value = 7750
$ws.Cells.Item(71, 13).Value = -4006

$ws.Cells.Item(122, 8).Value = 5026.6665
$ws.Cells.Item(122, 9).Value = 5026.6665
$ws.Cells.Item(122, 11).Value = 15079.9995
$ws.Cells.Item(122, 13).Value = -12629.9995

$ws.Cells.Item(125, 8).Value = 79775
$ws.Cells.Item(125, 10).Value = 79775
$ws.Cells.Item(125, 12).Value = 79775
$ws.Cells.Item(125, 14).Value = -89615

$ws.Cells.Item(126, 8).Value = 6908.3335
$ws.Cells.Item(126, 9).Value = 4900
$ws.Cells.Item(126, 10).Value = 8342.857
$ws.Cells.Item(126, 11).Value = 14700
$ws.Cells.Item(126, 12).Value = 25028.571
$ws.Cells.Item(126, 13).Value = -12230
$ws.Cells.Item(126, 14).Value = -29968.571

$ws.Cells.Item(132, 8).Value = 4900.394
$ws.Cells.Item(132, 9).Value = 4861.2173
$ws.Cells.Item(132, 10).Value = 4990.5
$ws.Cells.Item(132, 11).Value = 14583.6519
$ws.Cells.Item(132, 12).Value = 14971.5
$ws.Cells.Item(132, 13).Value = -12053.6519
$ws.Cells.Item(132, 14).Value = -20031.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3782.4211
$ws.Cells.Item(122, 9).Value = 2325.611
$ws.Cells.Item(122, 11).Value = 6976.833
$ws.Cells.Item(122, 13).Value = -4526.833

$ws.Cells.Item(132, 8).Value = 5025.5713
$ws.Cells.Item(132, 9).Value = 5295
$ws.Cells.Item(132, 10).Value = 4666.3335
$ws.Cells.Item(132, 11).Value = 5295.7
$ws.Cells.Item(132, 12).Value = 13999.0005
$ws.Cells.Item(132, 13).Value = -19059.0005

Write-Host "Updated 34 leve-profit rows across 8 sheets (210 cells)."
